$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric must be forced to Text format
# so Excel stores them as inline/shared strings instead of numbers,
# matching the original text-typed cells in the workbook.
$numericTextCells = @("D5","D6","D7","D8","D10","D11","D12","D16","D19","D20","D21","D22","D23","D24","D26","D29","D32","D33","D35","D37","D38","D39","D40","D41","D42","D44","D45","D46","D47","D48","D49","D51")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply new values
$ws.Range("D2").Value = "55.007.21"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").Value = "2.339.03"
$ws.Range("E3").Value = "  -4.32%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "474.76"
$ws.Range("E5").Value = "  -1.96%  "
$ws.Range("D6").Value = "144.97"
$ws.Range("E6").Value = "  +1.05%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "0.604"
$ws.Range("E8").Value = "  +19.81%  "
$ws.Range("D9").Value = "2.344.43"
$ws.Range("E9").Value = "  -4.35%  "
$ws.Range("D10").Value = "0.0955"
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("D11").Value = "5.43"
$ws.Range("E11").Value = "  -5.81%  "
$ws.Range("D12").Value = "0.323"
$ws.Range("E12").Value = "  -1.81%  "
$ws.Range("E13").Value = "  +1.24%  "
$ws.Range("D14").Value = "2.745.37"
$ws.Range("E14").Value = "  -4.51%  "
$ws.Range("D15").Value = "55.038.91"
$ws.Range("E15").Value = "  -1.53%  "
$ws.Range("D16").Value = "19.83"
$ws.Range("E16").Value = "  -4.69%  "
$ws.Range("E17").Value = "  -4.36%  "
$ws.Range("D18").Value = "2.343.01"
$ws.Range("E18").Value = "  -4.44%  "
$ws.Range("D19").Value = "4.56"
$ws.Range("E19").Value = "  +1.85%  "
$ws.Range("D20").Value = "313.04"
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("D21").Value = "9.49"
$ws.Range("E21").Value = "  -5.05%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "5.65"
$ws.Range("E23").Value = "  -2.33%  "
$ws.Range("D24").Value = "56.20"
$ws.Range("E24").Value = "  -3.51%  "
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").Value = "0.392"
$ws.Range("E26").Value = "  -3.63%  "
$ws.Range("E27").Value = "  -4.57%  "
$ws.Range("D28").Value = "2.447.30"
$ws.Range("E28").Value = "  -4.58%  "
$ws.Range("D29").Value = "7.02"
$ws.Range("E29").Value = "  -8.09%  "
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("D31").Value = "0.0₃0736"
$ws.Range("E31").Value = "  -4.60%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "18.10"
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "145.30"
$ws.Range("E33").Value = "  -1.62%  "
$ws.Range("E34").Value = "  -1.70%  "
$ws.Range("D35").Value = "5.07"
$ws.Range("E35").Value = "  -0.76%  "
$ws.Range("E36").Value = "  -4.05%  "
$ws.Range("D37").Value = "3.57"
$ws.Range("E37").Value = "  -2.66%  "
$ws.Range("D38").Value = "0.802"
$ws.Range("E38").Value = "  -4.95%  "
$ws.Range("D39").Value = "0.101"
$ws.Range("E39").Value = "  +10.40%  "
$ws.Range("D40").Value = "33.52"
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("D42").Value = "3.35"
$ws.Range("E42").Value = "  -3.94%  "
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("D44").Value = "0.575"
$ws.Range("E44").Value = "  -3.58%  "
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").Value = "10.14"
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "0.0513"
$ws.Range("E46").Value = "  -5.84%  "
$ws.Range("D47").Value = "247.75"
$ws.Range("E47").Value = "  -4.37%  "
$ws.Range("D48").Value = "0.0219"
$ws.Range("E48").Value = "  -2.90%  "
$ws.Range("D49").Value = "4.36"
$ws.Range("E49").Value = "  -6.34%  "
$ws.Range("D50").Value = "1.791.45"
$ws.Range("E50").Value = "  -3.50%  "
$ws.Range("D51").Value = "16.54"
$ws.Range("E51").Value = "  -4.49%  "

# Restore default (Normal) style on the cells we forced to Text format,
# so no stray style index is introduced while keeping the text type.
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).Style = "Normal"
}

Write-Output "done"
